$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.993.89"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "1.637.70"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'214.42"
$ws.Range("E5").Value = "  -1.01%  "

$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  -1.67%  "

$ws.Range("E9").Value = "  -2.41%  "

$ws.Range("E10").Value = "  -5.83%  "

$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "1.864.16"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").Value = "1.661.87"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("D14").Value = "'4.21"
$ws.Range("E14").Value = "  -2.01%  "

$ws.Range("D15").Value = "'0.530"
$ws.Range("E15").Value = "  -2.67%  "

$ws.Range("D16").Value = "25.991.00"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("E17").Value = "  -2.78%  "

$ws.Range("D18").Value = "'61.82"
$ws.Range("E18").Value = "  -2.26%  "

$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").Value = "'191.20"
$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").Value = "'4.25"
$ws.Range("E21").Value = "  -2.39%  "

$ws.Range("D22").Value = "'9.73"
$ws.Range("E22").Value = "  -2.22%  "

$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("E24").Value = "  +1.62%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'1.79"
$ws.Range("E25").Value = "  -1.00%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'143.72"
$ws.Range("E26").Value = "  -0.54%  "

$ws.Range("E27").Value = "  -0.25%  "

$ws.Range("D28").Value = "'6.85"
$ws.Range("E28").Value = "  -0.76%  "

$ws.Range("D29").Value = "'15.26"
$ws.Range("E29").Value = "  -1.79%  "

$ws.Range("E30").Value = "  -1.53%  "

$ws.Range("D31").Value = "'0.0483"
$ws.Range("E31").Value = "  -3.06%  "

$ws.Range("E32").Value = "  -3.06%  "

$ws.Range("E33").Value = "  -4.16%  "

$ws.Range("E34").Value = "  -2.35%  "

$ws.Range("E35").Value = "  -2.73%  "

$ws.Range("D36").Value = "1.135.52"
$ws.Range("E36").Value = "  +0.16%  "

$ws.Range("D37").Value = "'0.868"
$ws.Range("E37").Value = "  -4.24%  "

$ws.Range("E38").Value = "  -1.38%  "

$ws.Range("E39").Value = "  -3.83%  "

$ws.Range("E40").Value = "  -1.37%  "

$ws.Range("D41").Value = "'98.42"
$ws.Range("E41").Value = "  -1.11%  "

$ws.Range("E42").Value = "  -2.40%  "

$ws.Range("D43").Value = "1.773.22"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("E44").Value = "  -4.77%  "

$ws.Range("E45").Value = "  -0.88%  "

$ws.Range("D46").Value = "'55.24"
$ws.Range("E46").Value = "  -2.37%  "

$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("E48").Value = "  +2.73%  "

$ws.Range("E49").Value = "  -0.64%  "

$ws.Range("D50").Value = "'7.55"
$ws.Range("E50").Value = "  -2.24%  "

$ws.Range("E51").Value = "  -0.03%  "
